$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the EMIRO JOSE ORTIZ DURANGO rows (original rows 16-22, 7 rows).
#    This shifts the MARTHA LIGIA FONNEGRA GEORGE rows (originally 23-27) up to rows 16-20,
#    and the footer rows (originally 32-33) up to rows 25-26.
$ws.Range("B16:B22").EntireRow.Delete()

# 2) The remaining 5 "MARTHA" rows (now rows 16-20) were in descending period order
#    (1801, 1712, 1711, 1710, 1709). Re-order them ascending (1709, 1710, 1711, 1712, 1801)
#    to match the refreshed data. Valor Mora / Salario Basico stay 29600 / 840000 throughout.
$periods = @(1709, 1710, 1711, 1712, 1801)
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 5).Value = [string]$periods[$i]
}

# 3) Update the summary block with the refreshed totals.
$ws.Range("E11").Value = 148000          # VALOR MORA (was 402800)
$ws.Range("C13").Value = 1               # Cant. Trabajadores (was 2)
$ws.Range("F13").Value = 5               # Cant. Periodos (was 12)
